$wb = $excel.ActiveWorkbook

# The workbook contains a "展览" sheet and a "全部类型" sheet that mirrors it.
# Both need their "想去人数" (want-to-go count) column F updated for rows 2, 4, 5.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 356
    $ws.Range("F4").Value = 75
    $ws.Range("F5").Value = 293
}
